# Fixed a bug for respin stats
# The symbol rows (2-25) on the active sheet were shuffled/corrected.
# Apply the new row values (A:F) as given by the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(201, 9, 30, 15, 45, 30)
    3  = @(501, 9, 52, 30, 75, 45)
    4  = @(701, 3, 90, 45, 97, 15)
    5  = @(901, 16, 15, 45, 60, 60)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(301, 6, 45, 30, 60, 45)
    8  = @(601, 9, 60, 67, 60, 42)
    9  = @(801, 3, 67, 65, 52, 45)
    10 = @(1201, 2, 10, 10, 10, 10)
    11 = @(1202, 2, 10, 10, 10, 10)
    12 = @(101, 9, 30, 15, 60, 15)
    13 = @(401, 9, 48, 67, 75, 45)
    14 = @(1203, 3, 15, 15, 15, 15)
    15 = @(1001, 18, 30, 75, 60, 72)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(802, 0, 4, 5, 4, 0)
    18 = @(2, 0, 2, 2, 2, 2)
    19 = @(1, 0, 2, 2, 2, 2)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(3, 0, 3, 3, 3, 3)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
    24 = @(702, 0, 0, 0, 4, 0)
    25 = @(1002, 0, 0, 0, 0, 9)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
